$wb = $excel.ActiveWorkbook

# The two data sheets (same layout: column A = keys, column B = values)
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Overwrite every value cell (B2:B11) on both sheets with "test111".
# This also causes the now-unused shared strings (old credentials / sample
# values) to be dropped when the workbook is serialized.
for ($r = 2; $r -le 11; $r++) {
    $ws1.Cells.Item($r, 2).Value = "test111"
    $ws2.Cells.Item($r, 2).Value = "test111"
}

# Update the selection / active-sheet state:
#  - sheet2 loses its "tabSelected" flag and its selection becomes B2:B11
#  - sheet1 becomes the selected/active tab with selection B2:B11
$ws2.Range("B2:B11").Select()
$ws1.Range("B2:B11").Select()
